$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-case the header row labels (A1:C1) to title case.
$ws.Range("A1").Value = "Time"
$ws.Range("B1").Value = "Per"
$ws.Range("C1").Value = "Mode"

# Build the header format (bold, thin box border, centered/top aligned) on a
# scratch cell, then stamp it onto the header row with a single
# formats-only paste so the three header cells all share one new style
# (mirrors how Excel would fold one "apply header style" action into a
# single new cellXf instead of one per property tweak).
$tmpl = $ws.Range("F1")
$tmpl.Font.Bold = $true
$tmpl.Borders.LineStyle = 1
$tmpl.HorizontalAlignment = -4108
$tmpl.VerticalAlignment = -4160

$tmpl.Copy() | Out-Null
$hdr = $ws.Range("A1:C1")
$hdr.PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$tmpl.Clear() | Out-Null

# Move the active selection to E2.
$ws.Range("E2").Select() | Out-Null
